# "Add strength to attacks"
# Each hero's level-progression sheet computed the Strength stat (column D)
# at level 1 with a formula (=50-E2-F2-G2) that derived it from the other
# stats. Re-balance the starting Strength value for every hero by replacing
# that formula with a flat, hand-set number. Levels 2+ keep deriving from
# D2 via their existing formulas (D3 = D2*rate, D4:D11 = previous*rate),
# so they recalculate automatically once D2 changes.

$wb = $excel.ActiveWorkbook

# sheet name -> new flat Strength value for level 1 (cell D2)
$newD2 = @{
    "ulko_levels"    = 14
    "ferlin_levels"  = 8
    "phoebey_levels" = 12
    "bob_levels"     = 5
}

# sheet name -> new active-cell selection to restore after editing
$newSelection = @{
    "ulko_levels"    = "D16"
    "ferlin_levels"  = "D15"
    "phoebey_levels" = "D3"
    "bob_levels"     = "D15"
}

# Apply the Strength (D2) value change on every levels sheet.
foreach ($name in $newD2.Keys) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("D2").Value = $newD2[$name]
}

# Update the selection on the non-active sheets first so the originally
# active sheet (ulko_levels) ends up re-activated last and keeps its
# tabSelected state.
foreach ($name in $newSelection.Keys) {
    if ($name -ne "ulko_levels") {
        $ws = $wb.Worksheets.Item($name)
        [void]$ws.Range($newSelection[$name]).Select()
    }
}

$wsUlko = $wb.Worksheets.Item("ulko_levels")
[void]$wsUlko.Activate()
[void]$wsUlko.Range($newSelection["ulko_levels"]).Select()
